# Update column C (Förändrad) dates in rows 2-13 from 2023-10-13 (serial 45212)
# to 2023-10-22 (serial 45221) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45212) {
        $cell.Value2 = 45221
    }
}
